$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Color constants (RGB as used by Excel COM, decimal) ---
$yellow = 65535      # RGB(255,255,0)  -> fill used by style s=1
$green  = 5287936    # RGB(0,176,80)   -> fill used by style s=2 / s=5

# --- Row 47: update the "how" note for the overview-page requirement,
#     and re-color column A green (it was yellow before).
#     (Done before the "comments" edits below so new shared-string entries
#     come out in the same order as the target workbook: "partial" first.) ---
$ws.Range("B47").Value = "partial"
$ws.Range("A47").Interior.Color = $green

# --- Rows 14-16: move the "how" explanation from column B into column C,
#     and put "comments" in column B instead. Also highlight column A yellow. ---
$ws.Range("C14").Value = $ws.Range("B14").Value()
$ws.Range("B14").Value = "comments"
$ws.Range("A14").Interior.Color = $yellow

$ws.Range("C15").Value = $ws.Range("B15").Value()
$ws.Range("B15").Value = "comments"
$ws.Range("A15").Interior.Color = $yellow

$ws.Range("C16").Value = $ws.Range("B16").Value()
$ws.Range("B16").Value = "comments"
$ws.Range("A16").Interior.Color = $yellow

# --- Row 67: clear the stray "Migrasjonene..." bug note text, keep formatting. ---
$ws.Range("A67").Value = ""

# --- Color column A green for rows that previously had no fill. ---
$ws.Range("A32").Interior.Color = $green
$ws.Range("A33").Interior.Color = $green
$ws.Range("A34").Interior.Color = $green
$ws.Range("A36").Interior.Color = $green
$ws.Range("A37").Interior.Color = $green

$ws.Range("A50").Interior.Color = $green
$ws.Range("A51").Interior.Color = $green
$ws.Range("A52").Interior.Color = $green
$ws.Range("A53").Interior.Color = $green
$ws.Range("A54").Interior.Color = $green
$ws.Range("A55").Interior.Color = $green
$ws.Range("A56").Interior.Color = $green

$ws.Range("A61").Interior.Color = $green
$ws.Range("A62").Interior.Color = $green
$ws.Range("A63").Interior.Color = $green
$ws.Range("A64").Interior.Color = $green

# --- Update the selection / active cell shown when the sheet is opened. ---
$ws.Range("B17").Select()
